$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SYNoEVC")

# Row 1: shift the "start year" from 2021 to 2020 and drop all the other
# year columns (C1:AE1) - only the style is retained on the cells that had
# non-default styling.
$ws.Range("B1").Value = 2020
$ws.Range("C1:AE1").ClearContents()

# Row 2: point the lone remaining value at the new "2020" column of the
# Calculations sheet (row 8 instead of row 11) and clear out the rest of
# the year formulas (C2:AE2), keeping their style.
$ws.Range("B2").Formula = "=Calculations!B8"
$ws.Range("C2:AE2").ClearContents()

# Switch the active sheet/selection from "About" to "SYNoEVC", landing on
# B3 (matches the new tabSelected + selection in the saved view state).
$ws.Activate()
$ws.Range("B3").Select()
